# Generate Report for Handback
# For the f56291f8-515c-4d1a-b668-197c81cb9a58 file, the handback transform
# failed because the handback file name did not match the handoff file name.
# Update the Overview status and record the error detail on each language
# sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: row 3 is the f56291f8-515c-4d1a-b668-197c81cb9a58.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row 3 is the f56291f8-515c-4d1a-b668-197c81cb9a58 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "Handback file name: fvdb3ea0.x1q is different with handoff file name: f56291f8-515c-4d1a-b668-197c81cb9a58.aa80b538cd048d969c31194afc8293001f4faf6f.zh-cn."

# --- de-de sheet: row 3 is the f56291f8-515c-4d1a-b668-197c81cb9a58 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "Handback file name: fvdb3ea0.x1q is different with handoff file name: f56291f8-515c-4d1a-b668-197c81cb9a58.aa80b538cd048d969c31194afc8293001f4faf6f.de-de."
